$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44201
$ws.Range("I2").Value = "Segunda"
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 800
$ws.Range("L2").Value = 900
$ws.Range("M2").Value = 850
$ws.Range("P2").Value = 850

# Row 3
$ws.Range("D3").Value = 44229
$ws.Range("J3").Value = 760
$ws.Range("K3").Value = 550
$ws.Range("L3").Value = 600
$ws.Range("M3").Value = 575
$ws.Range("P3").Value = 575

# Row 4
$ws.Range("D4").Value = 44174
$ws.Range("J4").Value = 800
$ws.Range("K4").Value = 450
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = 475
$ws.Range("P4").Value = 475

# Row 5
$ws.Range("I5").Value = "Tercera"
$ws.Range("J5").Value = 1200
$ws.Range("K5").Value = 250
$ws.Range("L5").Value = 350
$ws.Range("M5").Value = 300
$ws.Range("P5").Value = 300

# Row 6
$ws.Range("D6").Value = 44245
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 850
$ws.Range("L6").Value = 900
$ws.Range("M6").Value = 875
$ws.Range("P6").Value = 875

# Row 7
$ws.Range("D7").Value = 44245
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 750
$ws.Range("L7").Value = 800
$ws.Range("M7").Value = 775
$ws.Range("P7").Value = 775

# Row 10
$ws.Range("D10").Value = 44267
$ws.Range("I10").Value = "Tercera"
$ws.Range("J10").Value = 400
$ws.Range("K10").Value = 500
$ws.Range("L10").Value = 600
$ws.Range("M10").Value = 550
$ws.Range("P10").Value = 550

# Row 11
$ws.Range("D11").Value = 44224
$ws.Range("J11").Value = 800
$ws.Range("K11").Value = 850
$ws.Range("L11").Value = 900
$ws.Range("M11").Value = 875
$ws.Range("P11").Value = 875

# Row 12
$ws.Range("D12").Value = 44210
$ws.Range("J12").Value = 900
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 700
$ws.Range("M12").Value = 650
$ws.Range("P12").Value = 650
